# Add a new weekly price record for "Vega Modelo de Temuco" (Maracuyá) as row 66,
# pushing the existing rows 66-98 down to 67-99.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("66:66").Insert()

$ws.Range("A66").Value = 10
$ws.Range("B66").Value = "Vega Modelo de Temuco"
$ws.Range("C66").Value = "La Araucanía"
$ws.Range("D66").Value = 45089
$ws.Range("E66").Value = 9
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100108
$ws.Range("H66").Value = "Tropicales y subtropicales"
$ws.Range("I66").Value = 100108003
$ws.Range("J66").Value = "Maracuyá"
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 50
$ws.Range("N66").Value = 50000
$ws.Range("O66").Value = 50000
$ws.Range("P66").Value = 50000
$ws.Range("Q66").Value = "$/caja 18 kilos"
$ws.Range("R66").Value = "Región de Arica y Parinacota"
$ws.Range("S66").Value = 2778
$ws.Range("T66").Value = 18
